$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (Vladimir Putin's original row) - this shifts rows 5 and 6 up to 4 and 5
$ws.Rows.Item(4).Delete()

# Append Vladimir Putin's updated record as the new row 6 (renumbered ID, new birthday,
# gender corrected, phone reused, address/picture unchanged)
$ws.Range("A6").Value = "2"
$ws.Range("B6").Value = "Vladimir"
$ws.Range("C6").Value = "Putin"
$ws.Range("D6").Value = 26851.6518209375
$ws.Range("E6").Value = "Female"
$ws.Range("F6").Value = "2399"
$ws.Range("G6").Value = "Russia"
$ws.Range("H6").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Putin.jpeg"
